$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16 - this shifts the existing rows 16-21 down to 17-22,
# matching the diff (old row16->17, 17->18, 18->19, 19->20, 20->21, 21->22).
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Cells.Item(16, 1).Value = 7
$ws.Cells.Item(16, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(16, 3).Value = 'Ñuble'
$ws.Cells.Item(16, 4).Value = 44489
$ws.Cells.Item(16, 5).Value = 16
$ws.Cells.Item(16, 6).Value = 100112013
$ws.Cells.Item(16, 7).Value = 'Alcachofa'
$ws.Cells.Item(16, 8).Value = 'Madrigal'
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 120
$ws.Cells.Item(16, 11).Value = 11000
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).Value = 11500
$ws.Cells.Item(16, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(16, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(16, 16).Value = 288
$ws.Cells.Item(16, 17).Value = 40
$ws.Cells.Item(16, 18).Value = 'Hortaliza'
